# This script shifts all "Starting time" (column A) and "End time" (column B)
# values in the Results sheet 4 hours earlier, e.g. "16:00" -> "12:00".
# The values are stored as plain text strings (e.g. "16:00"), not Excel time
# serials, so we parse/reformat them as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @(1, 2)) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Value2

        if ($null -ne $val -and $val -ne "") {
            $parts = $val.ToString().Split(":")
            if ($parts.Length -eq 2) {
                $hour = [int]$parts[0]
                $minute = [int]$parts[1]

                $hour = $hour - 4
                if ($hour -lt 0) {
                    $hour = $hour + 24
                }

                $newVal = "{0:D2}:{1:D2}" -f $hour, $minute
                $cell.Value2 = $newVal
            }
        }
    }
}
